$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells that contain numeric-looking values
# so Excel does not auto-convert them to numbers (preserving exact text formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update Price column values
$ws.Range("D2").Value = "69.974.28"
$ws.Range("D3").Value = "3.836.85"
$ws.Range("D4").Value = "0.998"
$ws.Range("D5").Value = "635.13"
$ws.Range("D6").Value = "166.70"
$ws.Range("D7").Value = "3.838.03"
$ws.Range("D11").Value = "0.454"
$ws.Range("D12").Value = "6.68"
$ws.Range("D14").Value = "36.02"
$ws.Range("D15").Value = "4.487.23"
$ws.Range("D16").Value = "3.857.79"
$ws.Range("D17").Value = "70.036.89"
$ws.Range("D18").Value = "18.14"
$ws.Range("D19").Value = "7.16"
$ws.Range("D21").Value = "471.62"
$ws.Range("D22").Value = "9.78"
$ws.Range("D23").Value = "0.710"
$ws.Range("D24").Value = "0.0000153"
$ws.Range("D25").Value = "84.03"
$ws.Range("D26").Value = "2.18"
$ws.Range("D27").Value = "12.07"
$ws.Range("D28").Value = "10.11"
$ws.Range("D30").Value = "3.991.33"
$ws.Range("D31").Value = "2.71"
$ws.Range("D32").Value = "7.33"
$ws.Range("D34").Value = "29.36"
$ws.Range("D35").Value = "3.788.09"
$ws.Range("D36").Value = "9.09"
$ws.Range("D37").Value = "0.998"
$ws.Range("D40").Value = "3.43"
$ws.Range("D41").Value = "5.93"
$ws.Range("D42").Value = "0.983"
$ws.Range("D43").Value = "0.999"
$ws.Range("D45").Value = "157.20"
$ws.Range("D46").Value = "0.302"
$ws.Range("D47").Value = "43.67"
$ws.Range("D48").Value = "1.95"
$ws.Range("D49").Value = "47.37"
$ws.Range("D51").Value = "8.47"

# Update other cells (Coin name, Link, Volume percentage)
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +5.64%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  +2.99%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("E39").Value = "  +8.91%  "
$ws.Range("E40").Value = "  +4.60%  "
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E48").Value = "  +3.45%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +4.64%  "
$ws.Range("E51").Value = "  +1.30%  "
